$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet1 ("Sheet1") - no data changes, just move the selection to C42.
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Activate()
$ws1.Range("C42").Select()

# ---------------------------------------------------------------------------
# Sheet2 ("Sheet2") - shift the dates forward a month (31 days), clear the
# last 6 rows out, auto-size column A, and move the view/selection.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()

for ($r = 2; $r -le 31; $r++) {
    $cell = $ws2.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 31
}

$ws2.Range("A32:B37").ClearContents()

$ws2.Columns.Item(1).AutoFit() | Out-Null

$ws2.Range("E32").Select()
$excel.ActiveWindow.ScrollRow = 13

# ---------------------------------------------------------------------------
# Sheet3 ("Sheet3") - shift the dates forward two months (61 days), replace
# row 33 with a marker label, clear the remaining trailing rows, auto-size
# column A, and move the view/selection. Do this sheet LAST so it stays the
# active / selected tab, matching the original workbook state.
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Activate()

for ($r = 2; $r -le 32; $r++) {
    $cell = $ws3.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 61
}

$ws3.Range("B33:B37").ClearContents()
$ws3.Range("A34:A37").ClearContents()
$ws3.Range("A33").Value = "^sppoky"

$ws3.Columns.Item(1).AutoFit() | Out-Null

$ws3.Range("A33").Select()
$excel.ActiveWindow.ScrollRow = 14

Write-Output "done"
